# Regenerate the "K" column (column G) values for the save_data sheet.
# These values were recomputed upstream (K instead of old Strike# calc,
# using regenerated std/mean and s_vals) and are written here directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(5, 7, 7, 3, 7, 9, 5, 5, 8, 12, 0, 9, 7, 5, 7, 4, 9, 5, 4, 3, 3, 7, 4, 4, 4, 4, 3, 2)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
